{"js": "// Extend the existing \"8. \" paragraph and append three new paragraphs\n// describing steps 9-11 (Routes/Controllers/ApiError, express-fileupload,\n// uuid), matching the target diff. The trailing run(s) with yellow\n// highlight are applied afterward via an in-paragraph search so the\n// run-splitting mirrors the author's formatting without depending on a\n// fragile insert+format-in-one-step sequence.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The last paragraph in the document is the \"8. \" list item that gets\n// extended with more text in this commit.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// --- Extend paragraph \"8. \" -------------------------------------------------\nlet tail = lastParagraph.getRange(\"End\");\ntail.insertText(\n  \"Creating models use sequelize lib after that in DB must be created tables \",\n  \"End\"\n);\nawait context.sync();\n\nlet hl = lastParagraph.search(\"sequelize lib\", { matchCase: true });\nhl.load(\"items\");\nawait context.sync();\nif (hl.items.length > 0) {\n  hl.items[0].font.highlightColor = \"Yellow\";\n  await context.sync();\n}\n\n// --- New paragraph \"9. ...\" --------------------------------------------------\nlet p9 = lastParagraph.insertParagraph(\n  \"9. Creating Routes, Controllers for each router, ApiError,  \",\n  \"After\"\n);\np9.paragraphFormat.leftIndent = 18; // 360 twips = 18pt (ind w:left=\"360\")\nawait context.sync();\n\n// --- New paragraph \"10. ...\" -------------------------------------------------\nlet p10 = p9.insertParagraph(\n  \"10. Install express-fileupload for const { img } = req.files in necessary Controler \",\n  \"After\"\n);\np10.paragraphFormat.leftIndent = 18;\nawait context.sync();\n\nlet hl10a = p10.search(\"express-fileupload\", { matchCase: true });\nhl10a.load(\"items\");\nlet hl10b = p10.search(\"const { img } = req.files\", { matchCase: true });\nhl10b.load(\"items\");\nawait context.sync();\nif (hl10a.items.length > 0) {\n  hl10a.items[0].font.highlightColor = \"Yellow\";\n}\nif (hl10b.items.length > 0) {\n  hl10b.items[0].font.highlightColor = \"Yellow\";\n}\nawait context.sync();\n\n// --- New paragraph \"11. ...\" -------------------------------------------------\nlet p11 = p10.insertParagraph(\n  \"11. Install uuid for generate unique id for uploading files \",\n  \"After\"\n);\np11.paragraphFormat.leftIndent = 18;\nawait context.sync();\n\nlet hl11 = p11.search(\"uuid\", { matchCase: true });\nhl11.load(\"items\");\nawait context.sync();\nif (hl11.items.length > 0) {\n  hl11.items[0].font.highlightColor = \"Yellow\";\n  await context.sync();\n}\n", "ps1": "# Extend the existing \"8. \" paragraph and append three new paragraphs\n# describing steps 9-11 (Routes/Controllers/ApiError, express-fileupload,\n# uuid), matching the target diff. Plain text is inserted first, then\n# Find.Execute is used to scope a highlight (wdYellow = 7) onto the\n# specific sub-strings that should carry the yellow highlight.\n\n$d = $word.ActiveDocument\n\n# The last paragraph in the document is the \"8. \" list item that gets\n# extended with more text in this commit.\n$count = $d.Paragraphs.Count\n$p8 = $d.Paragraphs.Item($count)\n\n# --- Extend paragraph \"8. \" --------------------------------------------\n$tail = $p8.Range\n$tail.Collapse(0)\n$tail.InsertAfter(\"Creating models use sequelize lib after that in DB must be created tables \")\n\n$p8 = $d.Paragraphs.Item($count)\n$find8 = $p8.Range\n$find8.Find.Text = \"sequelize lib\"\n$find8.Find.Execute() | Out-Null\nif ($find8.Find.Found) {\n    $find8.Font.HighlightColorIndex = 7\n}\n\n# --- New paragraph \"9. ...\" ---------------------------------------------\n$tail = $p8.Range\n$tail.Collapse(0)\n$tail.InsertParagraphAfter() | Out-Null\n$count = $d.Paragraphs.Count\n$p9 = $d.Paragraphs.Item($count)\n$p9.Range.Text = \"9. Creating Routes, Controllers for each router, ApiError,  \"\n$p9.Format.LeftIndent = 18\n\n# --- New paragraph \"10. ...\" ----------------------------------------------\n$tail = $p9.Range\n$tail.Collapse(0)\n$tail.InsertParagraphAfter() | Out-Null\n$count = $d.Paragraphs.Count\n$p10 = $d.Paragraphs.Item($count)\n$p10.Range.Text = \"10. Install express-fileupload for const { img } = req.files in necessary Controler \"\n$p10.Format.LeftIndent = 18\n\n$find10a = $p10.Range\n$find10a.Find.Text = \"express-fileupload\"\n$find10a.Find.Execute() | Out-Null\nif ($find10a.Find.Found) {\n    $find10a.Font.HighlightColorIndex = 7\n}\n\n$find10b = $p10.Range\n$find10b.Find.Text = \"const { img } = req.files\"\n$find10b.Find.Execute() | Out-Null\nif ($find10b.Find.Found) {\n    $find10b.Font.HighlightColorIndex = 7\n}\n\n# --- New paragraph \"11. ...\" ----------------------------------------------\n$tail = $p10.Range\n$tail.Collapse(0)\n$tail.InsertParagraphAfter() | Out-Null\n$count = $d.Paragraphs.Count\n$p11 = $d.Paragraphs.Item($count)\n$p11.Range.Text = \"11. Install uuid for generate unique id for uploading files \"\n$p11.Format.LeftIndent = 18\n\n$find11 = $p11.Range\n$find11.Find.Text = \"uuid\"\n$find11.Find.Execute() | Out-Null\nif ($find11.Find.Found) {\n    $find11.Font.HighlightColorIndex = 7\n}\n"}
